$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 67998
$ws.Range("I28").Value = 72847.86
$ws.Range("K28").Value = 72847.86
$ws.Range("M28").Value = -72362.86
$ws.Range("H70").Value = 84382.414
$ws.Range("J70").Value = 112233.555
$ws.Range("L70").Value = 336700.665
$ws.Range("N70").Value = -337240.665
$ws.Range("H73").Value = 84382.414
$ws.Range("J73").Value = 112233.555
$ws.Range("L73").Value = 336700.665
$ws.Range("N73").Value = -338572.665
$ws.Range("H81").Value = 33327.5
$ws.Range("J81").Value = 33327.5
$ws.Range("L81").Value = 33327.5
$ws.Range("N81").Value = -35323.5
$ws.Range("H84").Value = 33327.5
$ws.Range("J84").Value = 33327.5
$ws.Range("L84").Value = 99982.5
$ws.Range("N84").Value = -109966.5
$ws.Range("H86").Value = 4390121.5
$ws.Range("I86").Value = 2971
$ws.Range("J86").Value = 6583697
$ws.Range("K86").Value = 2971
$ws.Range("L86").Value = 6583697
$ws.Range("M86").Value = -1848
$ws.Range("N86").Value = -6585943
$ws.Range("H88").Value = 1823.4615
$ws.Range("I88").Value = 1802.5
$ws.Range("J88").Value = 1827.2727
$ws.Range("K88").Value = 1802.5
$ws.Range("L88").Value = 1827.2727
$ws.Range("M88").Value = -1396.5
$ws.Range("N88").Value = -2639.2727
$ws.Range("H89").Value = 4390121.5
$ws.Range("I89").Value = 2971
$ws.Range("J89").Value = 6583697
$ws.Range("K89").Value = 14855
$ws.Range("L89").Value = 32918485
$ws.Range("M89").Value = -9239
$ws.Range("N89").Value = -32929717
$ws.Range("H91").Value = 1823.4615
$ws.Range("I91").Value = 1802.5
$ws.Range("J91").Value = 1827.2727
$ws.Range("K91").Value = 1802.5
$ws.Range("L91").Value = 1827.2727
$ws.Range("M91").Value = -398.5
$ws.Range("N91").Value = -4635.2727
$ws.Range("H104").Value = 674.8333
$ws.Range("I104").Value = 674.8333
$ws.Range("K104").Value = 2024.4999
$ws.Range("M104").Value = -277.4999
$ws.Range("H112").Value = 3115.8696
$ws.Range("J112").Value = 3115.8696
$ws.Range("L112").Value = 9347.6088
$ws.Range("N112").Value = -11563.6088
$ws.Range("H124").Value = 75000
$ws.Range("J124").Value = 75000
$ws.Range("L124").Value = 75000
$ws.Range("N124").Value = -84820
$ws.Range("H125").Value = 18522454
$ws.Range("J125").Value = 27783004
$ws.Range("L125").Value = 250047036
$ws.Range("N125").Value = -250051956
$ws.Range("H135").Value = 590413.75
$ws.Range("I135").Value = 835066.25
$ws.Range("J135").Value = 3247.8
$ws.Range("K135").Value = 7515596.25
$ws.Range("L135").Value = 29230.2
$ws.Range("M135").Value = -7513061.25
$ws.Range("N135").Value = -34300.2
$ws.Range("H137").Value = 5177.0454
$ws.Range("I137").Value = 5042.619
$ws.Range("K137").Value = 15127.857
$ws.Range("M137").Value = -12577.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1991.579
$ws.Range("I102").Value = 2000.5
$ws.Range("J102").Value = 1944
$ws.Range("K102").Value = 2000.5
$ws.Range("L102").Value = 1944
$ws.Range("M102").Value = -378.5
$ws.Range("N102").Value = -5188
$ws.Range("H110").Value = 717455.5600000001
$ws.Range("I110").Value = 717455.5600000001
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 717455.5600000001
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -715410.5600000001
$ws.Range("H122").Value = 4945.1113
$ws.Range("I122").Value = 5366.9375
$ws.Range("K122").Value = 16100.8125
$ws.Range("M122").Value = -13650.8125
$ws.Range("H132").Value = 1941.85
$ws.Range("I132").Value = 1975.6666
$ws.Range("K132").Value = 5926.9998
$ws.Range("M132").Value = -3396.9998
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1064669
$ws.Range("J86").Value = 2600.5
$ws.Range("L86").Value = 2600.5
$ws.Range("N86").Value = -4846.5
$ws.Range("H89").Value = 1064669
$ws.Range("J89").Value = 2600.5
$ws.Range("L89").Value = 13002.5
$ws.Range("N89").Value = -24234.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 39500.355
$ws.Range("I31").Value = 1106
$ws.Range("J31").Value = 45899.418
$ws.Range("K31").Value = 1106
$ws.Range("L31").Value = 45899.418
$ws.Range("M31").Value = -811
$ws.Range("N31").Value = -46489.418
$ws.Range("H34").Value = 39500.355
$ws.Range("I34").Value = 1106
$ws.Range("J34").Value = 45899.418
$ws.Range("K34").Value = 1106
$ws.Range("L34").Value = 45899.418
$ws.Range("M34").Value = -904
$ws.Range("N34").Value = -46303.418

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1853499.4
$ws.Range("J113").Value = 1968.1428
$ws.Range("L113").Value = 5904.428400000001
$ws.Range("N113").Value = -10244.4284
$ws.Range("H127").Value = 1487.1111
$ws.Range("J127").Value = 1487.1111
$ws.Range("L127").Value = 4461.3333
$ws.Range("N127").Value = -14381.3333
$ws.Range("H131").Value = 7624974.5
$ws.Range("I131").Value = 66867770
$ws.Range("J131").Value = 29744.719
$ws.Range("K131").Value = 200603310
$ws.Range("L131").Value = 89234.15700000001
$ws.Range("M131").Value = -200598270
$ws.Range("N131").Value = -99314.15700000001
$ws.Range("H140").Value = 1909.2609
$ws.Range("I140").Value = 1677.8636
$ws.Range("K140").Value = 5033.5908
$ws.Range("M140").Value = 146.4092000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1004
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1004
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1004
$ws.Range("N5").Value = -1228
$ws.Range("H113").Value = 461302.2
$ws.Range("I113").Value = 910406.0600000001
$ws.Range("K113").Value = 910406.0600000001
$ws.Range("M113").Value = -908236.0600000001
$ws.Range("H122").Value = 4721.4
$ws.Range("I122").Value = 2982.8
$ws.Range("K122").Value = 8948.400000000001
$ws.Range("M122").Value = -6498.400000000001
$ws.Range("H123").Value = 54974.332
$ws.Range("J123").Value = 54974.332
$ws.Range("L123").Value = 54974.332
$ws.Range("N123").Value = -59874.332
$ws.Range("H132").Value = 53906.383
$ws.Range("I132").Value = 5391.0557
$ws.Range("J132").Value = 344998.34
$ws.Range("K132").Value = 16173.1671
$ws.Range("L132").Value = 1034995.02
$ws.Range("M132").Value = -13643.1671
$ws.Range("N132").Value = -1040055.02
$ws.Range("M5").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4631.4736
$ws.Range("I46").Value = 3964.1428
$ws.Range("K46").Value = 3964.1428
$ws.Range("M46").Value = -3776.1428
$ws.Range("H68").Value = 3199
$ws.Range("I68").Value = 2898.8
$ws.Range("J68").Value = 3949.5
$ws.Range("K68").Value = 2898.8
$ws.Range("L68").Value = 3949.5
$ws.Range("M68").Value = -2149.8
$ws.Range("N68").Value = -5447.5
$ws.Range("H71").Value = 3199
$ws.Range("I71").Value = 2898.8
$ws.Range("J71").Value = 3949.5
$ws.Range("K71").Value = 14494
$ws.Range("L71").Value = 19747.5
$ws.Range("M71").Value = -10750
$ws.Range("N71").Value = -27235.5
$ws.Range("H122").Value = 3030.0952
$ws.Range("I122").Value = 2703.2
$ws.Range("K122").Value = 8109.599999999999
$ws.Range("M122").Value = -5659.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8333.666999999999
$ws.Range("I62").Value = 10002
$ws.Range("K62").Value = 10002
$ws.Range("M62").Value = -9378
$ws.Range("H65").Value = 8333.666999999999
$ws.Range("I65").Value = 10002
$ws.Range("K65").Value = 50010
$ws.Range("M65").Value = -46890

